$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 342, shifting existing rows 342:365 down to 343:366
$ws.Rows("342:342").Insert()

# Populate the newly inserted row 342 with the new weekly price record
$ws.Cells.Item(342, 1).Value = 5
$ws.Cells.Item(342, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(342, 3).Value = "Maule"
$ws.Cells.Item(342, 4).Value = 44826
$ws.Cells.Item(342, 5).Value = 7
$ws.Cells.Item(342, 6).Value = 100112003
$ws.Cells.Item(342, 7).Value = "Ajo"
$ws.Cells.Item(342, 8).Value = "Chino"
$ws.Cells.Item(342, 9).Value = "Primera"
$ws.Cells.Item(342, 10).Value = 300
$ws.Cells.Item(342, 11).Value = 23000
$ws.Cells.Item(342, 12).Value = 23000
$ws.Cells.Item(342, 13).Value = 23000
$ws.Cells.Item(342, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(342, 15).Value = "China"
$ws.Cells.Item(342, 16).Value = 2300
$ws.Cells.Item(342, 17).Value = 10
$ws.Cells.Item(342, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date/number style as the other rows (column D)
$ws.Cells.Item(342, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
